# Rename the worksheet tab from "Specs" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Update D2 value from "--" to an em dash "—"
$ws.Range("D2").Value = "—"

# Delete row 3 (the "Sterility" row) entirely
$ws.Rows("3:3").Delete()
